$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4:D4").Value = "You need two saved frames"
